# Edit script for LOQ4206.docx
#
# The edit reshuffles several text blocks between paragraphs (content is
# moved/rewritten in place; paragraph count, order, styles and run
# formatting - e.g. bold labels, italic EN text - stay the same).
#
# Paragraph indices below are the *original* (before-edit) 1-based
# Word `Paragraphs` indices; since no paragraphs are inserted or removed,
# they remain stable targets throughout the whole script.

$d = $word.ActiveDocument

function Replace-InParagraph {
    param(
        [int]$ParaIndex,
        [string]$OldText,
        [string]$NewText
    )
    $rng = $d.Paragraphs.Item($ParaIndex).Range
    $found = $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)
    if (-not $found) {
        throw "Replace-InParagraph: text not found in paragraph $ParaIndex"
    }
}
# NOTE: this minimal PowerShell host does not bind named (`-Param value`)
# arguments reliably, so Replace-InParagraph is always invoked positionally
# below: Replace-InParagraph <ParaIndex> <OldText> <NewText>

$nl = [char]11   # manual line break (w:br) inside a Word Range.Text assignment

# --- Paragraph 6 (PT "Objetivos" body, plain run) --------------------------
# becomes the PT "Programa resumido" summary text
$d.Paragraphs.Item(6).Range.Text = "Introdução a Pesquisa Operacional, Programação Linear, Método Simplex, Introdução aos Grafos e à Otimização em Rede, Estudo de Casos em Programação Linear, Introdução a Teoria das Filas,"

# --- Paragraph 7 (EN "Objetivos" body, italic run) --------------------------
# becomes the EN "Programa resumido" summary text
$d.Paragraphs.Item(7).Range.Text = "Introduction to Operational Research, Linear Programming, Simplex Method, Introduction to Graphs and Network Optimization, Case Study in Linear Programming, Introduction to Queue Theory."

# --- Paragraph 9 (ListBullet, "Docente(s)" value) ---------------------------
# becomes the old PT "Objetivos" body text
$d.Paragraphs.Item(9).Range.Text = "Compreender a Pesquisa Operacional como ciência aplicada à Engenharia de Produção. Proporcionar conhecimento dos problemas típicos encontrados em Engenharia de Produção. Analisar, Modelar e solucionar os problemas por meio da Pesquisa Operacional."

# --- Paragraph 11 (PT "Programa resumido" body, plain run) ------------------
# becomes the full PT "Programa" outline (many lines separated by manual breaks)
$programaPT = @(
    "1. Introdução a Pesquisa Operacional",
    "1.1. Conceitos de Pesquisa Operacional;",
    "1.2. Modelagem;",
    "1.3. Estrutura dos Modelos Matemáticos;",
    "1.4. Técnicas matemáticas em Pesquisa Operacional;",
    "1.2. Fases de Um Estudo em Pesquisa Operacional",
    "2. Programação Linear",
    "2.1. Definição",
    "2.2. Formulação de Modelos",
    "2.3. Resolução Gráfica;",
    "3. Método Simplex",
    "3.1. Desenvolvimento do Método Simplex;",
    "3.2. Procedimento do Método Simplex;",
    "4. Introdução aos Grafos e à Otimização em Rede",
    "4.1. Conceitos Básicos em Teoria dos Grafos",
    "4.2. Problemas de Fluxo Máximo;",
    "4.3. Problemas de Caminho Mínimo",
    "5. Estudo de Casos em Programação Linear",
    "5.1. Modelo de Transporte Simples",
    "5.2. Modelo da Designação.",
    "6. Introdução a Teoria das Filas",
    "6.1. Conceitos da Teoria das Filas",
    "6.2. Modelos Markovianos"
) -join $nl
$d.Paragraphs.Item(11).Range.Text = $programaPT

# --- Paragraph 12 (EN "Programa resumido" body, italic run) -----------------
# becomes the old EN "Objetivos" body text
$d.Paragraphs.Item(12).Range.Text = "Understand Operational Research as a science applied to Industrial Engineering. Provide knowledge of the typical problems encountered in Industrial Engineering. Analyze, model and solve problems through Operational Research."

# --- Paragraph 14 (PT "Programa" outline body, plain run) -------------------
# becomes the "Avaliação" / "Método:" value text
$d.Paragraphs.Item(14).Range.Text = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# --- Paragraph 17 (ListBullet, "Avaliação" Método/Critério/Norma block) -----
# Only the three non-bold "value" runs change text; the bold "Método: " /
# "Critério: " / "Norma de recuperação: " label runs are untouched.
# Do the replacements from the last value run to the first, because each
# new value is equal to a value used further down the paragraph (old text
# of a later run); replacing in document order would make a later Find
# re-match text just inserted by an earlier step.
$bibNew = @(
    "1. HILLIER, F.S., LIEBERMAN, G.J., Introdução à Pesquisa Operacional, 8ªed., Editora McGraw-Hill, 2006.",
    "2. LACHTERMACHER, G., Pesquisa Operacional na Tomada de Decisão (modelagem em Excel), 4ª ed., Editora Campus, 2009.",
    "3. ANDERSON, D.R., SWEENEY, D.J. e WILLIAMS, T.A., An Introduction to Management Science 9ª ed., South-Western College Publishing, 2000.",
    "4. PIZZOLATO, N. D. e GANDOLPHO, A. A. Técnicas de Otimização, LTC Editora, 2009.",
    "5. TAHA, H. A ., Pesquisa Operacional, 8ª ed., Pearson/Prentice Hall, 2008."
) -join $nl

Replace-InParagraph 17 `
    "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada." `
    $bibNew

Replace-InParagraph 17 `
    "NF≥ 5,0." `
    "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."

Replace-InParagraph 17 `
    "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n." `
    "NF≥ 5,0."

# --- Paragraph 19 (PT "Bibliografia" body, plain run) ------------------------
# becomes the "Docente(s)" value text
$d.Paragraphs.Item(19).Range.Text = "5840917 - Fabricio Maciel Gomes"

Write-Output "Edit complete."
